$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K14").Value = 2.1
$ws.Range("L14").Value = 2.62
$ws.Range("M14").Value = 1.01
$ws.Range("N14").Value = 8.1
$ws.Range("P14").Value = 2.95
$ws.Range("T14").Value = 2.55
$ws.Range("V14").Value = 1.85
$ws.Range("W14").Value = 9.75
$ws.Range("AD14").Value = 6.7
$ws.Range("AK14").Value = 16.5
$ws.Range("AR14").Value = 150
$ws.Range("AY14").Value = 20
